$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 388, shifting the existing rows (388-426) down to (389-427).
$ws.Rows(388).Insert()

# Populate the newly inserted row with the new weekly price-report record.
$ws.Cells.Item(388, 1).Value  = 10
$ws.Cells.Item(388, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(388, 3).Value  = "La Araucanía"
$ws.Cells.Item(388, 4).Value  = 45194
$ws.Cells.Item(388, 5).Value  = 9
$ws.Cells.Item(388, 6).Value  = 100112039
$ws.Cells.Item(388, 7).Value  = "Ciboulette"
$ws.Cells.Item(388, 8).Value  = "Sin especificar"
$ws.Cells.Item(388, 9).Value  = "Primera"
$ws.Cells.Item(388, 10).Value = 40
$ws.Cells.Item(388, 11).Value = 7000
$ws.Cells.Item(388, 12).Value = 7000
$ws.Cells.Item(388, 13).Value = 7000
$ws.Cells.Item(388, 14).Value = "$/docena de atados"
$ws.Cells.Item(388, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(388, 16).Value = 2333
$ws.Cells.Item(388, 17).Value = 3
$ws.Cells.Item(388, 18).Value = "Hortaliza"
